# The presentation originally carries two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (clrScheme "Office") - unused, only
#                            wired to the Notes Master relationship.
#   ppt/theme/theme2.xml -> "Integral" (clrScheme "Red Violet") - the theme
#                            actually driving the Slide Master / slides.
#
# The target edit swaps the content of the two theme parts. The Slide
# Master's live theme (persisted as ppt/theme/theme2.xml) needs to end up
# holding the "Office" color scheme (the font scheme and format scheme are
# already byte-identical between the two theme parts, so only the 12-slot
# color scheme actually changes visibly).
#
# Apply this through the real PowerPoint object model: the editable surface
# for a theme's colors is SlideMaster.Theme.ThemeColorScheme, whose 12 items
# follow the standard MsoThemeColorSchemeIndex order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1..accent6, 11 hlink, 12 folHlink

function Get-VbaRgb {
    param([int]$r, [int]$g, [int]$b)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$scheme = $theme.ThemeColorScheme

# Target "Office" color scheme (the colors currently sitting, unused, in
# ppt/theme/theme1.xml) expressed as R,G,B triples.
$officeColors = @(
    @(0x00, 0x00, 0x00),  # 1  dk1
    @(0xFF, 0xFF, 0xFF),  # 2  lt1
    @(0x44, 0x54, 0x6A),  # 3  dk2
    @(0xE7, 0xE6, 0xE6),  # 4  lt2
    @(0x5B, 0x9B, 0xD5),  # 5  accent1
    @(0xED, 0x7D, 0x31),  # 6  accent2
    @(0xA5, 0xA5, 0xA5),  # 7  accent3
    @(0xFF, 0xC0, 0x00),  # 8  accent4
    @(0x44, 0x72, 0xC4),  # 9  accent5
    @(0x70, 0xAD, 0x47),  # 10 accent6
    @(0x05, 0x63, 0xC1),  # 11 hlink
    @(0x95, 0x4F, 0x72)   # 12 folHlink
)

for ($i = 1; $i -le $scheme.Count; $i++) {
    $rgb = $officeColors[$i - 1]
    $scheme.Item($i).RGB = Get-VbaRgb $rgb[0] $rgb[1] $rgb[2]
}
